$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.859557867050171
$ws.Range("B1").Value = 2.775059223175049
$ws.Range("C1").Value = 2.963214874267578
$ws.Range("D1").Value = 3.5876624584198
$ws.Range("E1").Value = 1.652899980545044
